$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") rows 2-22: update date serial value from 45212 to 45221
for ($row = 2; $row -le 22; $row++) {
    $ws.Cells.Item($row, 3).Value2 = 45221
}
